$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add row 19, mirroring the pattern/content of row 9 (same posting date,
# company, country, material doc, etc.) but with a reversal quantity
# (O19 = -0.5) whose std-price total is still 67.5.

# A19: copy the date style from A9 (keeps the existing numFmtId=164 style
# index instead of minting a new cellXf), then set its own value.
$ws.Range("A9").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Cells.Item(19, 1).Value = 44928

# B19:N19 (F19 stays blank, same as row 9): copy values+types from row 9
# so the shared-string text cells come through as t="s" without picking up
# any unwanted formatting.
$ws.Range("B9:N9").Copy()
$ws.Range("B19").PasteSpecial(-4104)

# O19 / P19: plain numeric cells.
$ws.Cells.Item(19, 15).Value = -0.5
$ws.Cells.Item(19, 16).Value = 67.5

# Match the author's final selection state.
$ws.Range("V14").Select()
